# Manual data update: append the next day's COVID-19 statistics row to the
# "Covid-19 podatki" sheet / "Tabela1" table (row 82, date 31/5/2020).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Covid-19 podatki")
$lo = $ws.ListObjects.Item("Tabela1")

# Grow the table by one row - this keeps the table ref / autofilter / dimension
# in sync automatically (mirrors what Excel does when you type into the row
# directly below a table).
$newRow = $lo.ListRows.Add()

# Carry the formatting of the preceding "normal" data row down onto the new
# row (row 81 itself carries a one-off formatting quirk in column B, so pull
# from row 80 which uses the standard per-column styling instead).
$ws.Range("A80:J80").Copy($ws.Range("A82:J82")) | Out-Null

$rng = $newRow.Range
$rng.Item(1, 1).Value = 43982
$rng.Item(1, 2).Value = 79039
$rng.Item(1, 3).Value = 246
$rng.Item(1, 4).Value = 1473
$rng.Item(1, 5).Value = 0
$rng.Item(1, 6).Value = 5
$rng.Item(1, 7).Value = 1
$rng.Item(1, 8).Value = 0
$rng.Item(1, 9).Value = 109
$rng.Item(1, 10).Value = 1

$ws.Range("A82:J82").Select() | Out-Null
